$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 9899.5
$ws.Range("I32").Value = 9801
$ws.Range("J32").Value = 9998
$ws.Range("K32").Value = 9801
$ws.Range("L32").Value = 9998
$ws.Range("M32").Value = -9475
$ws.Range("N32").Value = -10650

$ws.Range("H98").Value = 2806180.2
$ws.Range("I98").Value = 3368640
$ws.Range("K98").Value = 3368640
$ws.Range("M98").Value = -3367142

$ws.Range("H112").Value = 6972745.5
$ws.Range("J112").Value = 7746950.5
$ws.Range("L112").Value = 23240851.5
$ws.Range("N112").Value = -23243067.5

$ws.Range("H115").Value = 67347650
$ws.Range("I115").Value = 67347650
$ws.Range("K115").Value = 202042950
$ws.Range("M115").Value = -202041383

$ws.Range("H122").Value = 2806180.2
$ws.Range("I122").Value = 3368640
$ws.Range("K122").Value = 10105920
$ws.Range("M122").Value = -10103470

$ws.Range("H132").Value = 4022.1724
$ws.Range("I132").Value = 1887
$ws.Range("J132").Value = 22527
$ws.Range("K132").Value = 5661
$ws.Range("L132").Value = 67581
$ws.Range("M132").Value = -3131
$ws.Range("N132").Value = -72641

$ws.Range("H137").Value = 1298809.9
$ws.Range("I137").Value = 2249274.8
$ws.Range("J137").Value = 2721.5454
$ws.Range("K137").Value = 6747824.399999999
$ws.Range("L137").Value = 8164.6362
$ws.Range("M137").Value = -6745274.399999999
$ws.Range("N137").Value = -13264.6362

$ws.Range("H138").Value = 1655.45
$ws.Range("I138").Value = 871
$ws.Range("J138").Value = 2024.6029
$ws.Range("K138").Value = 2613
$ws.Range("L138").Value = 6073.8087
$ws.Range("M138").Value = 2527
$ws.Range("N138").Value = -16353.8087

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7685.7
$ws.Range("I32").Value = 6320.7812
$ws.Range("K32").Value = 6320.7812
$ws.Range("M32").Value = -6033.7812

$ws.Range("H74").Value = 40872.08
$ws.Range("I74").Value = 44661.02
$ws.Range("J74").Value = 5256
$ws.Range("K74").Value = 44661.02
$ws.Range("L74").Value = 5256
$ws.Range("M74").Value = -43787.02
$ws.Range("N74").Value = -7004

$ws.Range("H77").Value = 40872.08
$ws.Range("I77").Value = 44661.02
$ws.Range("J77").Value = 5256
$ws.Range("K77").Value = 223305.1
$ws.Range("L77").Value = 26280
$ws.Range("M77").Value = -218937.1
$ws.Range("N77").Value = -35016

$ws.Range("H122").Value = 4142.4443
$ws.Range("I122").Value = 2852.9644
$ws.Range("J122").Value = 8655.625
$ws.Range("K122").Value = 8558.893199999999
$ws.Range("L122").Value = 25966.875
$ws.Range("M122").Value = -6108.893199999999
$ws.Range("N122").Value = -30866.875

$ws.Range("H132").Value = 2223.6272
$ws.Range("I132").Value = 2272.875
$ws.Range("J132").Value = 1304.3334
$ws.Range("K132").Value = 6818.625
$ws.Range("L132").Value = 3913.0002
$ws.Range("M132").Value = -4288.625
$ws.Range("N132").Value = -8973.0002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3473.5454
$ws.Range("I99").Value = 2929.9
$ws.Range("J99").Value = 4638.5
$ws.Range("K99").Value = 2929.9
$ws.Range("L99").Value = 4638.5
$ws.Range("M99").Value = -1431.9
$ws.Range("N99").Value = -7634.5

$ws.Range("H105").Value = 1250.7742
$ws.Range("I105").Value = 1276.8518
$ws.Range("K105").Value = 1276.8518
$ws.Range("M105").Value = 470.1482000000001

$ws.Range("H134").Value = 3527.6365
$ws.Range("I134").Value = 2985.5
$ws.Range("J134").Value = 4973.3335
$ws.Range("K134").Value = 8956.5
$ws.Range("L134").Value = 14920.0005
$ws.Range("M134").Value = -6421.5
$ws.Range("N134").Value = -19990.0005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 32606760
$ws.Range("I31").Value = 501613.6
$ws.Range("J31").Value = 90979750
$ws.Range("K31").Value = 501613.6
$ws.Range("L31").Value = 90979750
$ws.Range("M31").Value = -501318.6
$ws.Range("N31").Value = -90980340

$ws.Range("H34").Value = 32606760
$ws.Range("I34").Value = 501613.6
$ws.Range("J34").Value = 90979750
$ws.Range("K34").Value = 501613.6
$ws.Range("L34").Value = 90979750
$ws.Range("M34").Value = -501411.6
$ws.Range("N34").Value = -90980154

$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("M50").ClearContents()

$ws.Range("H58").Value = 2359.5173
$ws.Range("I58").Value = 2272.3572
$ws.Range("K58").Value = 2272.3572
$ws.Range("M58").Value = -2069.3572

$ws.Range("H80").Value = 48998
$ws.Range("J80").Value = 48998
$ws.Range("L80").Value = 48998
$ws.Range("N80").Value = -51244

$ws.Range("H83").Value = 48998
$ws.Range("J83").Value = 48998
$ws.Range("L83").Value = 146994
$ws.Range("N83").Value = -158226

$ws.Range("H87").Value = 80000
$ws.Range("J87").Value = 80000
$ws.Range("L87").Value = 80000
$ws.Range("N87").Value = -82372

$ws.Range("H90").Value = 80000
$ws.Range("J90").Value = 80000
$ws.Range("L90").Value = 240000
$ws.Range("N90").Value = -251856

$ws.Range("H105").Value = 4460.778
$ws.Range("I105").Value = 1253.2142
$ws.Range("J105").Value = 6501.9546
$ws.Range("K105").Value = 1253.2142
$ws.Range("L105").Value = 6501.9546
$ws.Range("M105").Value = 493.7858000000001
$ws.Range("N105").Value = -9995.954600000001

$ws.Range("H132").Value = 2270.6216
$ws.Range("I132").Value = 959.25806
$ws.Range("K132").Value = 2877.77418
$ws.Range("M132").Value = -347.7741799999999

$ws.Range("H134").Value = 3646.725
$ws.Range("I134").Value = 3825.8157
$ws.Range("K134").Value = 11477.4471
$ws.Range("M134").Value = -8942.447100000001

$ws.Range("H136").Value = 2359.5173
$ws.Range("I136").Value = 2272.3572
$ws.Range("K136").Value = 6817.071599999999
$ws.Range("M136").Value = -4267.071599999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 6495102.5
$ws.Range("I131").Value = 50000788
$ws.Range("J131").Value = 1716.194
$ws.Range("K131").Value = 150002364
$ws.Range("L131").Value = 5148.582
$ws.Range("M131").Value = -149997324
$ws.Range("N131").Value = -15228.582

$ws.Range("H137").Value = 18521080
$ws.Range("I137").Value = 1525.4445
$ws.Range("J137").Value = 37040636
$ws.Range("K137").Value = 4576.333500000001
$ws.Range("L137").Value = 111121908
$ws.Range("M137").Value = 523.6664999999994
$ws.Range("N137").Value = -111132108

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 5435.9
$ws.Range("I113").Value = 3959.25
$ws.Range("J113").Value = 7650.875
$ws.Range("K113").Value = 3959.25
$ws.Range("L113").Value = 7650.875
$ws.Range("M113").Value = -1789.25
$ws.Range("N113").Value = -11990.875

$ws.Range("H126").Value = 4676.1113
$ws.Range("I126").Value = 4806.4287
$ws.Range("J126").Value = 4220
$ws.Range("K126").Value = 14419.2861
$ws.Range("L126").Value = 12660
$ws.Range("M126").Value = -11949.2861
$ws.Range("N126").Value = -17600

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 22957
$ws.Range("I43").Value = 22957
$ws.Range("J43").Value = 22957
$ws.Range("K43").Value = 22957
$ws.Range("L43").Value = 22957
$ws.Range("M43").Value = -22764
$ws.Range("N43").Value = -23343

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1633.5625
$ws.Range("I126").Value = 1675.8
$ws.Range("K126").Value = 5027.4
$ws.Range("M126").Value = -2557.4

$ws.Range("H132").Value = 15386736
$ws.Range("I132").Value = 19232686
$ws.Range("J132").Value = 2936.6924
$ws.Range("K132").Value = 57698058
$ws.Range("L132").Value = 8810.0772
$ws.Range("M132").Value = -57695528
$ws.Range("N132").Value = -13870.0772

$ws.Range("H136").Value = 199926.06
$ws.Range("I136").Value = 273796
$ws.Range("J136").Value = 4698.357
$ws.Range("K136").Value = 821388
$ws.Range("L136").Value = 14095.071
$ws.Range("M136").Value = -818838
$ws.Range("N136").Value = -19195.071
